# Fruta / hortaliza, semanal
# The sheet holds one row per market observation (rows 2-32). This edit
# re-shuffles the weekly observations: for each row, the date (D), Volumen
# (J), Precio minimo/maximo/promedio (K/L/M) and Precio $/Kg (P) values are
# replaced with those from another row in the original data (a pure
# permutation of those six columns across rows 2-32; every other column -
# market, region, category, unit, origin, etc. - stays put).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# new row -> source row (in the ORIGINAL/before data) that supplies the
# new D/J/K/L/M/P values for that row.
$mapping = @{
    2 = 7;   3 = 16;  4 = 32;  5 = 31;  6 = 28;  7 = 13;  8 = 3;
    9 = 11;  10 = 18; 11 = 5;  12 = 2;  13 = 19; 14 = 10; 15 = 21;
    16 = 8;  17 = 17; 18 = 24; 19 = 20; 20 = 12; 21 = 4;  22 = 29;
    23 = 27; 24 = 30; 25 = 15; 26 = 14; 27 = 25; 28 = 6;  29 = 23;
    30 = 9;  31 = 26; 32 = 22
}

$firstRow = 2
$lastRow = 32

# Snapshot the original values first so the re-shuffle reads are not
# clobbered by earlier writes in the same pass.
$snapD = @{}
$snapJ = @{}
$snapK = @{}
$snapL = @{}
$snapM = @{}
$snapP = @{}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $snapD[$r] = $ws.Cells.Item($r, 4).Value2
    $snapJ[$r] = $ws.Cells.Item($r, 10).Value2
    $snapK[$r] = $ws.Cells.Item($r, 11).Value2
    $snapL[$r] = $ws.Cells.Item($r, 12).Value2
    $snapM[$r] = $ws.Cells.Item($r, 13).Value2
    $snapP[$r] = $ws.Cells.Item($r, 16).Value2
}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $src = $mapping[$r]
    $ws.Cells.Item($r, 4).Value2 = $snapD[$src]
    $ws.Cells.Item($r, 10).Value2 = $snapJ[$src]
    $ws.Cells.Item($r, 11).Value2 = $snapK[$src]
    $ws.Cells.Item($r, 12).Value2 = $snapL[$src]
    $ws.Cells.Item($r, 13).Value2 = $snapM[$src]
    $ws.Cells.Item($r, 16).Value2 = $snapP[$src]
}

Write-Output "Done re-shuffling rows $firstRow..$lastRow"
